$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.374.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.802.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("E6").Value = '  +4.15%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '36.20'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.98%  '
$ws.Range("E9").Value = '  +2.18%  '
$ws.Range("E10").Value = '  +0.69%  '
$ws.Range("E11").Value = '  +2.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.062.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.798.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("E15").Value = '  +1.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '34.339.97'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.69%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.33%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.67%  '
$ws.Range("E27").Value = '  +2.49%  '
$ws.Range("E28").Value = '  +2.36%  '
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.14%  '
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.398.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.672'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("E37").Value = '  -3.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.07'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("E40").Value = '  +10.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.965'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '82.62'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0507'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.02'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.962.20'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '104.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").Value = '0.0₆0128'
$ws.Range("E51").Value = '  +0.41%  '
